# Apply the 2024-09-12 crypto price/volume refresh described in the commit
# message ("Updated cryptos list ... with GitHub Actions"). Rows 13/14,
# 26/27 and 46/47 also swap rank order (Avalanche <-> wstETH, ICP <->
# BSC-USD, Injective <-> Hedera) in addition to their value refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.709.19"
$ws.Range("E2").Value = "  +3.14%  "
# Row 3
$ws.Range("D3").Value = "2.332.57"
$ws.Range("E3").Value = "  +1.52%  "
# Row 4
$ws.Range("E4").Value = "  -0.01%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.09"
$ws.Range("E5").Value = "  +5.30%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.96"
$ws.Range("E6").Value = "  +3.89%  "
# Row 7
$ws.Range("E7").Value = "  +0.22%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("E8").Value = "  +6.75%  "
# Row 9
$ws.Range("E9").Value = "  +3.18%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.45"
$ws.Range("E10").Value = "  +4.72%  "
# Row 11
$ws.Range("E11").Value = "  +0.03%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("E12").Value = "  +7.26%  "
# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.71"
$ws.Range("E13").Value = "  +2.64%  "
# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.750.62"
$ws.Range("E14").Value = "  +1.61%  "
# Row 15
$ws.Range("D15").Value = "57.678.51"
$ws.Range("E15").Value = "  +3.22%  "
# Row 16
$ws.Range("E16").Value = "  +1.96%  "
# Row 17
$ws.Range("D17").Value = "2.303.65"
$ws.Range("E17").Value = "  -0.30%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.64"
$ws.Range("E18").Value = "  +3.47%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "334.34"
$ws.Range("E19").Value = "  +2.51%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("E20").Value = "  +3.55%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.68"
$ws.Range("E21").Value = "  +0.35%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.33%  "
# Row 23
$ws.Range("E23").Value = "  +1.04%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.01"
$ws.Range("E24").Value = "  +2.12%  "
# Row 25
$ws.Range("E25").Value = "  +3.00%  "
# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.51"
$ws.Range("E26").Value = "  +0.27%  "
# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.27%  "
# Row 28
$ws.Range("E28").Value = "  +7.00%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.76"
$ws.Range("E29").Value = "  +5.15%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.39"
$ws.Range("E30").Value = "  +2.04%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0731"
$ws.Range("E31").Value = "  +3.73%  "
# Row 32
$ws.Range("E32").Value = "  +1.90%  "
# Row 33
$ws.Range("E33").Value = "  +18.63%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.46"
$ws.Range("E34").Value = "  +1.75%  "
# Row 35
$ws.Range("E35").Value = "  +0.02%  "
# Row 36
$ws.Range("E36").Value = "  +0.29%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.17"
$ws.Range("E37").Value = "  +8.47%  "
# Row 38
$ws.Range("E38").Value = "  +2.65%  "
# Row 39
$ws.Range("E39").Value = "  +4.58%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.08"
$ws.Range("E40").Value = "  +2.09%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.07"
$ws.Range("E41").Value = "  -0.50%  "
# Row 42
$ws.Range("E42").Value = "  +1.06%  "
# Row 43
$ws.Range("E43").Value = "  +2.71%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "285.63"
$ws.Range("E44").Value = "  +2.29%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0936"
$ws.Range("E45").Value = "  +1.57%  "
# Row 46
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.04"
$ws.Range("E46").Value = "  +6.68%  "
# Row 47
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0503"
$ws.Range("E47").Value = "  +2.18%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.561"
$ws.Range("E48").Value = "  +1.59%  "
# Row 49
$ws.Range("E49").Value = "  +1.92%  "
# Row 50
$ws.Range("E50").Value = "  +2.36%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.47"
$ws.Range("E51").Value = "  +2.36%  "
